$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark that currently sits at the end
#    of the Introduction paragraph (after the final ".").
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Locate the sentence in the Methods paragraph that needs to be
#    split and extended with the new GAUL reference.
# ------------------------------------------------------------------
$sentence = $d.Content
$sentence.Find.Execute("was used to extract the resulting area under the disease progress curves by state and were classified according to relative risk.")
$sentenceStart = $sentence.Start
$sentenceEnd = $sentence.End

$marker = "by state"
$splitOffset = $sentence.Text.IndexOf($marker) + $marker.Length
$splitPoint = $sentenceStart + $splitOffset

$newText = " using FAO’s Global Administrative Unit Layers (GAUL)"

# ------------------------------------------------------------------
# 3. Insert the new text right after "...by state".
# ------------------------------------------------------------------
$insertionPoint = $d.Range($splitPoint, $splitPoint)
$insertionPoint.InsertAfter($newText)

# ------------------------------------------------------------------
# 4. Any text edit above can cause the engine to re-merge runs that
#    share identical (empty) formatting, so the run breaks are
#    (re-)established last, using temporary bookmarks that are added
#    and immediately removed again - this splits the run without
#    leaving any stray run-formatting behind.
#
#    Break #1: between "...(Hijmans 2015) " and "was used...".
# ------------------------------------------------------------------
$leftBreakRange = $d.Range($sentenceStart, $sentenceStart)
$d.Bookmarks.Add("ZZZTempSplit0", $leftBreakRange)
$d.Bookmarks("ZZZTempSplit0").Delete()

# Break #2: between "...by state" and " using FAO's...".
$tempBookmarkRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("ZZZTempSplit", $tempBookmarkRange)
$d.Bookmarks("ZZZTempSplit").Delete()

# ------------------------------------------------------------------
# 5. Re-insert the "_GoBack" bookmark right after the newly added
#    text (this both marks the last-edit location, like Word does,
#    and forces the run break between the new text and the
#    remainder of the sentence).
# ------------------------------------------------------------------
$goBackPoint = $splitPoint + $newText.Length
$goBackRange = $d.Range($goBackPoint, $goBackPoint)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# ------------------------------------------------------------------
# 6. Break #3: between "...relative risk." and " Results were
#    compared...", which the earlier text edit also merged away.
# ------------------------------------------------------------------
$rightBreakPoint = $sentenceEnd + $newText.Length
$rightBreakRange = $d.Range($rightBreakPoint, $rightBreakPoint)
$d.Bookmarks.Add("ZZZTempSplit2", $rightBreakRange)
$d.Bookmarks("ZZZTempSplit2").Delete()
